$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run Info")

# Copy the date/time number format (style) from A13 down through the new rows
# so the new column-A cells keep the same style index as the existing ones.
$ws.Range("A13").Copy()
$ws.Range("A14:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @{ r=14; A=42945.580775462964;  B="rcp60"; C=1; D=1000; E=4; F=10.197927302969147; G=9.615384615384615;  H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=15; A=42945.733703703707;  B="rcp60"; C=1; D=1000; E=4; F=19.929067068743567; G=9.615384615384615;  H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=16; A=42945.743715277778;  B="rcp60"; C=1; D=1000; E=4; F=5.336427200732647;  G=9.615384615384615;  H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=17; A=42945.744398148148;  B="rcp60"; C=1; D=1000; E=4; F=5.0318168305917208; G=9.615384615384615;  H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=18; A=42945.752025462964;  B="rcp60"; C=1; D=1000; E=8; F=12.080413886628522; G=13.01775147928994;  H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=19; A=42945.752662037034;  B="rcp60"; C=1; D=1000; E=8; F=7.3886378930385428; G=13.01775147928994;  H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=20; A=42945.759548611109;  B="rcp85"; C=1; D=1000; E=8; F=6.8727726762635646; G=11.834319526627219; H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 },
    @{ r=21; A=42945.759872685187;  B="rcp85"; C=0; D=1000; E=8; F=6.5432937176231087; G=11.538461538461538; H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.2549999999999999 },
    @{ r=22; A=42945.922500000001;  B="rcp85"; C=0; D=1000; E=8; F=24.113477714523054; G=11.538461538461538; H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.2549999999999999 },
    @{ r=23; A=42945.923935185187;  B="rcp85"; C=1; D=1000; E=4; F=5.4166161350635358; G=8.6538461538461533; H=0.3; I=0.1; J=4; K=4; L=2; M=0.36; N=1.5; O=0.46; P=4.7156000000000002 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
}

$ws.Range("A23:P23").Select()
